# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the
# c160a3b9-... file row (row 3) on both the "zh-cn" and "de-de" sheets,
# reflecting a freshly regenerated handback report.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 3 (c160a3b9 file) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-11 16:38:08"
$wsZhCn.Range("G3").Value = "2016-01-11 16:39:48"

# --- de-de sheet: row 3 (c160a3b9 file) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-11 16:38:35"
$wsDeDe.Range("G3").Value = "2016-01-11 16:40:27"
